# KichbanDOAN.docx -- "cap nhat lai kich ban" (re-work the change-log text)
#
# Summary of the content edits (see the unified diff of word/document.xml):
#   P2 (1.0): drop the trailing ", chuc nang tim kim" clause
#   P3 (1.1): "trang quan ly Admin giao dien" -> ", chuc nang tim kiem"
#   P4 (1.2): drop the trailing ", quan ly mac hang va sang pham." clause
#   P5 (2.0): reworded to "... trang quan ly Admin giao dien, quan ly mat
#             hang va san pham" and the old hoa-don/Admin sentence (with its
#             _GoBack bookmark) is pushed into a brand-new P6 ("Cap nhat 2.1"),
#             which also gains a trailing ", ".
#
# A small helper trick is used throughout to force Word to materialise a new
# <w:r> run boundary at a given point without altering the visible
# formatting: toggle Bold on then back off across the sub-range. Because the
# range shares its right edge with an already-settled boundary (end of the
# paragraph/clause being split), only one new edge -- at the range's start --
# is actually introduced.

$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Paragraph 2 ("Cap nhat 1.0 ...")
#   "...danh muc san pham, chuc nang tim kim." ->
#     run1 "...danh muc san pha"  +  run2 "m."
# ------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute(", chức năng tìm kím", $true, $false, $false, $false, $false, `
                 $true, 1, $false, "", 0) | Out-Null
$r.Text = ""

$r = $d.Content
$r.Find.Execute("sản phẩm.", $true, $false, $false, $false, $false, `
                 $true, 1, $false, "", 0) | Out-Null
$r.MoveStart(1, 7) | Out-Null   # isolate the trailing "m."
$r.Bold = 1
$r.Bold = 0

# ------------------------------------------------------------------
# Paragraph 3 ("Cap nhat 1.1 ...")
#   "...gio hang, trang quan ly Admin giao dien." ->
#     run1 "...gio hang" + run2 "," + run3 " " + run4 "chuc nang tim kiem."
# ------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("trang quản lý Admin giao diện.", $true, $false, $false, `
                 $false, $false, $true, 1, $false, "", 0) | Out-Null
$r.Text = "chức năng tìm kiếm."

$r = $d.Content
$r.Find.Execute("giỏ hàng, chức năng tìm kiếm.", $true, $false, $false, `
                 $false, $false, $true, 1, $false, "", 0) | Out-Null
$base = $r.Start
$end = $r.End

$b1 = $d.Range($base + 8, $end)     # after "giỏ hàng"
$b1.Bold = 1
$b1.Bold = 0

$b2 = $d.Range($base + 9, $end)     # after ","
$b2.Bold = 1
$b2.Bold = 0

$b3 = $d.Range($base + 10, $end)    # after the following " "
$b3.Bold = 1
$b3.Bold = 0

# ------------------------------------------------------------------
# Paragraph 4 ("Cap nhat 1.2 ...")
#   drop ", quan ly mac hang va sang pham." and split the remainder:
#     run1 "Cap nhat 1.2 Cap nha"  +  run2 "t tinh nang thanh toan"
# ------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute(", quản lý mặc hàng và sảng phẩm.", $true, $false, $false, `
                 $false, $false, $true, 1, $false, "", 0) | Out-Null
$r.Text = ""

$r = $d.Content
$r.Find.Execute("Cập nhật 1.2 Cập nhật tính năng thanh toán", $true, $false, `
                 $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$base = $r.Start
$end = $r.End

$b1 = $d.Range($base + 20, $end)    # "Cập nhật 1.2 Cập nhậ" | "t tính ..."
$b1.Bold = 1
$b1.Bold = 0

# ------------------------------------------------------------------
# Paragraph 5 ("Cap nhat 2.0 ...") -> split into paragraph 5 + new paragraph 6
# ------------------------------------------------------------------

# Split right after "Cap nhat 2.0": the remainder of the old sentence
# (hoa-don/Admin text, bookmark included) becomes a brand-new paragraph.
$r = $d.Content
$r.Find.Execute("Cập nhật 2.0", $true, $false, $false, $false, $false, `
                 $true, 1, $false, "", 0) | Out-Null
$r.Collapse(0) | Out-Null
$r.InsertParagraphAfter()

# New paragraph 6 currently reads " Cap nhat tinh nang ... boi Admin)".
# Prefix it with "Cap nhat 2.1" (and restore the Times New Roman font, since
# InsertBefore at this fresh paragraph/run boundary falls back to the
# document default font otherwise).
$r = $d.Content
$r.Find.Execute(" Cập nhật tính năng quản lý hóa đơn", $true, $false, `
                 $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r.Collapse(1) | Out-Null
$r.InsertBefore("Cập nhật 2.1")

$r = $d.Content
$r.Find.Execute("Cập nhật 2.1", $true, $false, $false, $false, $false, `
                 $true, 1, $false, "", 0) | Out-Null
$r.Font.NameAscii = "Times New Roman"
$r.Font.NameOther = "Times New Roman"
$r.Font.NameBi = "Times New Roman"

# Append ", " after the closing ")" at the end of paragraph 6 (same font fix).
$p6 = $d.Paragraphs(6)
$r = $p6.Range
$r.MoveEnd(1, -1) | Out-Null   # exclude the paragraph mark
$r.Collapse(0) | Out-Null
$r.InsertAfter(", ")

$r = $d.Content
$r.Find.Execute("bởi Admin), ", $true, $false, $false, $false, $false, `
                 $true, 1, $false, "", 0) | Out-Null
$r.MoveStart(1, $r.End - $r.Start - 2) | Out-Null   # isolate the trailing ", "
$r.Font.NameAscii = "Times New Roman"
$r.Font.NameOther = "Times New Roman"
$r.Font.NameBi = "Times New Roman"

# Paragraph 5 now holds just "Cap nhat 2.0"; append the reworded sentence as
# three further runs:
#   " trang quan ly Admin giao dien," + " " + "quan ly mat hang va san pham"
$p5 = $d.Paragraphs(5)
$r = $p5.Range
$r.MoveEnd(1, -1) | Out-Null   # exclude the paragraph mark
$r.Collapse(0) | Out-Null
$r.InsertAfter(" trang quản lý Admin giao diện, quản lý mặt hàng và sản phẩm")

$r = $d.Content
$r.Find.Execute(" trang quản lý Admin giao diện, quản lý mặt hàng và sản phẩm", `
                 $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$base = $r.Start
$end = $r.End

$b0 = $d.Range($base, $end)          # "Cập nhật 2.0" | " trang ..."
$b0.Bold = 1
$b0.Bold = 0

$b1 = $d.Range($base + 31, $end)     # "...giao diện," | " " | "quản lý ..."
$b1.Bold = 1
$b1.Bold = 0

$b2 = $d.Range($base + 32, $end)
$b2.Bold = 1
$b2.Bold = 0

Write-Output "edit complete"
